$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Midterm 1 grade update ---
# Fill in newly-graded Midterm 1 scores (column F) for students who
# previously had no grade recorded.
$ws.Range("F4").Formula  = "=40/60"
$ws.Range("F11").Formula = "=55/60"
$ws.Range("F12").Formula = "=54/60"
$ws.Range("F13").Formula = "=52/60"
$ws.Range("F15").Formula = "=50/60"
$ws.Range("F20").Formula = "=(60-14)/60"
$ws.Range("F25").Formula = "=51/60"
$ws.Range("F34").Formula = "=40/60"

# Re-confirm / re-enter grades for students whose Midterm 1 score was
# already present (values unchanged, but re-entered as part of this
# grading pass).
$ws.Range("F3").Formula  = "=51/60"
$ws.Range("F5").Formula  = "=61/60"
$ws.Range("F6").Formula  = "=51/60"
$ws.Range("F9").Formula  = "=40/60"
$ws.Range("F14").Formula = "=41/60"
$ws.Range("F28").Formula = "=58/60"
$ws.Range("F32").Formula = "=42/60"
$ws.Range("F35").Formula = "=40/60"

# Row 15 grew slightly taller once its Midterm 1 score was entered.
$ws.Rows(15).RowHeight = 15

# Leave the selection on F12, matching where editing finished.
$ws.Range("F12").Select() | Out-Null
